$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number (45186 = 2023-09-17)
# that was bumped to 45188 (2023-09-19) for every data row (rows 2-157).
$ws.Range("C2:C157").Value = 45188
